# Auto-generated edit script: updates the 'F' column (想去人数 / interested-count)
# values across the four worksheets, per the source diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 823
$ws.Range("F3").Value = 14753
$ws.Range("F5").Value = 1665
$ws.Range("F6").Value = 509
$ws.Range("F7").Value = 2137
$ws.Range("F8").Value = 1308
$ws.Range("F9").Value = 1992
$ws.Range("F10").Value = 951
$ws.Range("F12").Value = 2361
$ws.Range("F13").Value = 623
$ws.Range("F14").Value = 839
$ws.Range("F15").Value = 3678
$ws.Range("F17").Value = 347
$ws.Range("F18").Value = 2734
$ws.Range("F19").Value = 702
$ws.Range("F20").Value = 132
$ws.Range("F22").Value = 1938
$ws.Range("F23").Value = 1141
$ws.Range("F24").Value = 1659
$ws.Range("F26").Value = 180
$ws.Range("F27").Value = 7657
$ws.Range("F28").Value = 5290
$ws.Range("F29").Value = 334
$ws.Range("F31").Value = 728
$ws.Range("F32").Value = 734
$ws.Range("F33").Value = 3417
$ws.Range("F35").Value = 931
$ws.Range("F36").Value = 365
$ws.Range("F37").Value = 158
$ws.Range("F38").Value = 127
$ws.Range("F39").Value = 4520
$ws.Range("F40").Value = 749
$ws.Range("F41").Value = 36
$ws.Range("F42").Value = 354

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F13").Value = 19
$ws.Range("F15").Value = 103
$ws.Range("F17").Value = 111
$ws.Range("F18").Value = 125
$ws.Range("F19").Value = 61
$ws.Range("F25").Value = 24

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 8069
$ws.Range("F3").Value = 326
$ws.Range("F4").Value = 1151

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 8069
$ws.Range("F3").Value = 823
$ws.Range("F4").Value = 326
$ws.Range("F5").Value = 1151
$ws.Range("F6").Value = 14753
$ws.Range("F9").Value = 1665
$ws.Range("F10").Value = 509
$ws.Range("F11").Value = 1308
$ws.Range("F12").Value = 1992
$ws.Range("F15").Value = 623
$ws.Range("F17").Value = 3678
$ws.Range("F18").Value = 347
$ws.Range("F19").Value = 2734
$ws.Range("F20").Value = 702
$ws.Range("F21").Value = 132
$ws.Range("F23").Value = 1938
$ws.Range("F27").Value = 19
$ws.Range("F29").Value = 1659
$ws.Range("F30").Value = 103
$ws.Range("F32").Value = 180
$ws.Range("F33").Value = 7658
$ws.Range("F34").Value = 5290
$ws.Range("F35").Value = 334
$ws.Range("F36").Value = 728
$ws.Range("F37").Value = 3417
$ws.Range("F39").Value = 931
$ws.Range("F40").Value = 365
$ws.Range("F42").Value = 127
$ws.Range("F43").Value = 4520
$ws.Range("F44").Value = 749
$ws.Range("F45").Value = 36
$ws.Range("F46").Value = 354
$ws.Range("F48").Value = 24
